$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 25
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 53
$ws.Range("F18").Value = 37
$ws.Range("G18").Value = 90
